$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.6695229441184551
$ws.Cells.Item(2, 4).Value = 0.0783750512527206
$ws.Cells.Item(2, 5).Value = 0.06141686888949849
$ws.Cells.Item(2, 6).Value = 1.8225937690151
$ws.Cells.Item(2, 7).Value = 0.002509639840078791
$ws.Cells.Item(2, 9).Value = 1.314337654316148
$ws.Cells.Item(2, 11).Value = 1.078039710455187
$ws.Cells.Item(2, 12).Value = 0.3048057343004587
$ws.Cells.Item(2, 14).Value = 2.384885438280094
$ws.Cells.Item(3, 2).Value = 0.648343806870173
$ws.Cells.Item(3, 4).Value = 0.07896966201160449
$ws.Cells.Item(3, 5).Value = 0.06002831827511734
$ws.Cells.Item(3, 6).Value = 1.795894908131046
$ws.Cells.Item(3, 7).Value = 0.002514287915409635
$ws.Cells.Item(3, 9).Value = 1.316302648305587
$ws.Cells.Item(3, 11).Value = 0.9849463595662087
$ws.Cells.Item(3, 12).Value = 0.2901764589762905
$ws.Cells.Item(3, 14).Value = 2.400058022474838
$ws.Cells.Item(4, 2).Value = 0.6356953306500372
$ws.Cells.Item(4, 4).Value = 0.07934965744824751
$ws.Cells.Item(4, 5).Value = 0.05916249308003074
$ws.Cells.Item(4, 6).Value = 1.780555789770219
$ws.Cells.Item(4, 7).Value = 0.002517293043700939
$ws.Cells.Item(4, 9).Value = 1.318029098011991
$ws.Cells.Item(4, 11).Value = 0.9283347762168432
$ws.Cells.Item(4, 12).Value = 0.2813799228976279
$ws.Cells.Item(4, 14).Value = 2.410068766200837
$ws.Cells.Item(5, 2).Value = 0.6306306491133
$ws.Cells.Item(5, 4).Value = 0.07950825326515698
$ws.Cells.Item(5, 5).Value = 0.05880626684290213
$ws.Cells.Item(5, 6).Value = 1.774569258071097
$ws.Cells.Item(5, 7).Value = 0.002518555801429298
$ws.Cells.Item(5, 9).Value = 1.318863307216866
$ws.Cells.Item(5, 11).Value = 0.9054018969729043
$ws.Cells.Item(5, 12).Value = 0.2778418774154829
$ws.Cells.Item(5, 14).Value = 2.414322702786066
$ws.Cells.Item(6, 2).Value = 0.6297950864016855
$ws.Cells.Item(6, 4).Value = 0.0795348140708878
$ws.Cells.Item(6, 5).Value = 0.05874690882044753
$ws.Cells.Item(6, 6).Value = 1.77359113863136
$ws.Cells.Item(6, 7).Value = 0.002518767788749491
$ws.Cells.Item(6, 9).Value = 1.319009716834557
$ws.Cells.Item(6, 11).Value = 0.9016021426400869
$ws.Cells.Item(6, 12).Value = 0.2772571992795463
$ws.Cells.Item(6, 14).Value = 2.415039597491351
$ws.Cells.Item(7, 2).Value = 0.6356266631498499
$ws.Cells.Item(7, 4).Value = 0.07935178116931318
$ws.Cells.Item(7, 5).Value = 0.05915770271879062
$ws.Cells.Item(7, 6).Value = 1.780473984271126
$ws.Cells.Item(7, 7).Value = 0.002517309918962431
$ws.Cells.Item(7, 9).Value = 1.318039819489385
$ws.Cells.Item(7, 11).Value = 0.9280249428327636
$ws.Cells.Item(7, 12).Value = 0.2813320190107618
$ws.Cells.Item(7, 14).Value = 2.410125430124808
$ws.Cells.Item(8, 2).Value = 0.6621467177052125
$ws.Cells.Item(8, 4).Value = 0.07857697856183066
$ws.Cells.Item(8, 5).Value = 0.06094080773689292
$ws.Cells.Item(8, 6).Value = 1.813168743299599
$ws.Cells.Item(8, 7).Value = 0.002511211192152696
$ws.Cells.Item(8, 9).Value = 1.31490722378711
$ws.Cells.Item(8, 11).Value = 1.045826711548642
$ws.Cells.Item(8, 12).Value = 0.2997228922546782
$ws.Cells.Item(8, 14).Value = 2.389972627860274
$ws.Cells.Item(9, 2).Value = 0.7169663153418355
$ws.Cells.Item(9, 4).Value = 0.0771759276359596
$ws.Cells.Item(9, 5).Value = 0.06433535032801529
$ws.Cells.Item(9, 6).Value = 1.885687654080968
$ws.Cells.Item(9, 7).Value = 0.002500445517988934
$ws.Cells.Item(9, 9).Value = 1.312894158389661
$ws.Cells.Item(9, 11).Value = 1.281243909968794
$ws.Cells.Item(9, 12).Value = 0.3372705772967777
$ws.Cells.Item(9, 14).Value = 2.355974874821285
$ws.Cells.Item(10, 2).Value = 0.75895270107182
$ws.Cells.Item(10, 4).Value = 0.07621880888587906
$ws.Cells.Item(10, 5).Value = 0.06677142093856681
$ws.Cells.Item(10, 6).Value = 1.944155976801369
$ws.Cells.Item(10, 7).Value = 0.00249325579554283
$ws.Cells.Item(10, 9).Value = 1.313940704080956
$ws.Cells.Item(10, 11).Value = 1.456993138813004
$ws.Cells.Item(10, 12).Value = 0.3657757861977302
$ws.Cells.Item(10, 14).Value = 2.334376087867795
$ws.Cells.Item(11, 2).Value = 0.7784241332295778
$ws.Cells.Item(11, 4).Value = 0.0757990988583872
$ws.Cells.Item(11, 5).Value = 0.06786804152660331
$ws.Cells.Item(11, 6).Value = 1.971896086988806
$ws.Cells.Item(11, 7).Value = 0.00249013959563935
$ws.Cells.Item(11, 9).Value = 1.314967039049662
$ws.Cells.Item(11, 11).Value = 1.537574497681931
$ws.Cells.Item(11, 12).Value = 0.3789464795558217
$ws.Cells.Item(11, 14).Value = 2.325286951908197
$ws.Cells.Item(12, 2).Value = 0.7858507306826539
$ws.Cells.Item(12, 4).Value = 0.07564242539537158
$ws.Cells.Item(12, 5).Value = 0.06828171639558356
$ws.Cells.Item(12, 6).Value = 1.982565887360465
$ws.Cells.Item(12, 7).Value = 0.002488981650056333
$ws.Cells.Item(12, 9).Value = 1.315434940128526
$ws.Cells.Item(12, 11).Value = 1.568180952628609
$ws.Cells.Item(12, 12).Value = 0.3839633544061201
$ws.Cells.Item(12, 14).Value = 2.321951271887812
$ws.Cells.Item(13, 2).Value = 0.7842489183926205
$ws.Cells.Item(13, 4).Value = 0.0756760671488852
$ws.Cells.Item(13, 5).Value = 0.06819269375789538
$ws.Cells.Item(13, 6).Value = 1.980260591815238
$ws.Cells.Item(13, 7).Value = 0.002489230053882575
$ws.Cells.Item(13, 9).Value = 1.315330642573365
$ws.Cells.Item(13, 11).Value = 1.561585192940584
$ws.Cells.Item(13, 12).Value = 0.3828815697790162
$ws.Cells.Item(13, 14).Value = 2.32266494302641
$ws.Cells.Item(14, 2).Value = 0.7790340593068379
$ws.Cells.Item(14, 4).Value = 0.07578616387803905
$ws.Cells.Item(14, 5).Value = 0.06790210629486992
$ws.Cells.Item(14, 6).Value = 1.972770581258203
$ws.Cells.Item(14, 7).Value = 0.002490043888761135
$ws.Cells.Item(14, 9).Value = 1.315003944538766
$ws.Cells.Item(14, 11).Value = 1.540090657871929
$ws.Cells.Item(14, 12).Value = 0.3793586303886798
$ws.Cells.Item(14, 14).Value = 2.325010393764629
$ws.Cells.Item(15, 2).Value = 0.7758467269399887
$ws.Cells.Item(15, 4).Value = 0.07585389600633707
$ws.Cells.Item(15, 5).Value = 0.06772390793510574
$ws.Cells.Item(15, 6).Value = 1.968204277943727
$ws.Cells.Item(15, 7).Value = 0.002490545259729022
$ws.Cells.Item(15, 9).Value = 1.31481415701704
$ws.Cells.Item(15, 11).Value = 1.526936655392262
$ws.Cells.Item(15, 12).Value = 0.3772045640701833
$ws.Cells.Item(15, 14).Value = 2.326460886783821
$ws.Cells.Item(16, 2).Value = 0.7576876536615771
$ws.Cells.Item(16, 4).Value = 0.0762465540871391
$ws.Cells.Item(16, 5).Value = 0.06669952761259523
$ws.Cells.Item(16, 6).Value = 1.942366165433327
$ws.Cells.Item(16, 7).Value = 0.002493462544442751
$ws.Cells.Item(16, 9).Value = 1.313884713787061
$ws.Cells.Item(16, 11).Value = 1.451739775110639
$ws.Cells.Item(16, 12).Value = 0.3649191588777683
$ws.Cells.Item(16, 14).Value = 2.334984925687081
$ws.Cells.Item(17, 2).Value = 0.7466426515769911
$ws.Cells.Item(17, 4).Value = 0.07649145914987354
$ws.Cells.Item(17, 5).Value = 0.06606819268340836
$ws.Cells.Item(17, 6).Value = 1.926808613582324
$ws.Cells.Item(17, 7).Value = 0.00249529167649694
$ws.Cells.Item(17, 9).Value = 1.313455551727287
$ws.Cells.Item(17, 11).Value = 1.405771640017463
$ws.Cells.Item(17, 12).Value = 0.3574346969431588
$ws.Cells.Item(17, 14).Value = 2.340402936887543
$ws.Cells.Item(18, 2).Value = 0.7403248513349183
$ws.Cells.Item(18, 4).Value = 0.07663379834354522
$ws.Cells.Item(18, 5).Value = 0.06570397313928922
$ws.Cells.Item(18, 6).Value = 1.917967810196075
$ws.Cells.Item(18, 7).Value = 0.002496358288708555
$ws.Cells.Item(18, 9).Value = 1.313260494196662
$ws.Cells.Item(18, 11).Value = 1.379391416503381
$ws.Cells.Item(18, 12).Value = 0.3531489804659884
$ws.Cells.Item(18, 14).Value = 2.343588515414012
$ws.Cells.Item(19, 2).Value = 0.7381917704190357
$ws.Cells.Item(19, 4).Value = 0.07668224533306311
$ws.Cells.Item(19, 5).Value = 0.06558046490086511
$ws.Cells.Item(19, 6).Value = 1.914992899509357
$ws.Cells.Item(19, 7).Value = 0.002496721926766754
$ws.Cells.Item(19, 9).Value = 1.313203341468999
$ws.Cells.Item(19, 11).Value = 1.370469699142518
$ws.Cells.Item(19, 12).Value = 0.3517011946257043
$ws.Cells.Item(19, 14).Value = 2.344678987555625
$ws.Cells.Item(20, 2).Value = 0.7478147910976816
$ws.Cells.Item(20, 4).Value = 0.07646523575962672
$ws.Cells.Item(20, 5).Value = 0.06613551191053979
$ws.Cells.Item(20, 6).Value = 1.928453610230619
$ws.Cells.Item(20, 7).Value = 0.002495095458016709
$ws.Cells.Item(20, 9).Value = 1.313495876323927
$ws.Cells.Item(20, 11).Value = 1.410658867277732
$ws.Cells.Item(20, 12).Value = 0.3582294478498227
$ws.Cells.Item(20, 14).Value = 2.339819007786815
$ws.Cells.Item(21, 2).Value = 0.7805643482736571
$ws.Cells.Item(21, 4).Value = 0.07575376437304637
$ws.Cells.Item(21, 5).Value = 0.06798750150178279
$ws.Cells.Item(21, 6).Value = 1.974966087419261
$ws.Cells.Item(21, 7).Value = 0.002489804246857105
$ws.Cells.Item(21, 9).Value = 1.315097751837271
$ws.Cells.Item(21, 11).Value = 1.546401616654066
$ws.Cells.Item(21, 12).Value = 0.3803926034372864
$ws.Cells.Item(21, 14).Value = 2.324318593915507
$ws.Cells.Item(22, 2).Value = 0.802277988144624
$ws.Cells.Item(22, 4).Value = 0.07530196197607708
$ws.Cells.Item(22, 5).Value = 0.0691886405757689
$ws.Cells.Item(22, 6).Value = 2.006328160499606
$ws.Cells.Item(22, 7).Value = 0.002486474847677421
$ws.Cells.Item(22, 9).Value = 1.316606661028359
$ws.Cells.Item(22, 11).Value = 1.635654390638422
$ws.Cells.Item(22, 12).Value = 0.3950490561368554
$ws.Cells.Item(22, 14).Value = 2.314807212336547
$ws.Cells.Item(23, 2).Value = 0.7906607344553152
$ws.Cells.Item(23, 4).Value = 0.07554188900713488
$ws.Cells.Item(23, 5).Value = 0.06854839214622821
$ws.Cells.Item(23, 6).Value = 1.989501165729934
$ws.Cells.Item(23, 7).Value = 0.002488240071936142
$ws.Cells.Item(23, 9).Value = 1.315759013862944
$ws.Cells.Item(23, 11).Value = 1.58796898473048
$ws.Cells.Item(23, 12).Value = 0.3872108839599946
$ws.Cells.Item(23, 14).Value = 2.319826871755083
$ws.Cells.Item(24, 2).Value = 0.747284766757673
$ws.Cells.Item(24, 4).Value = 0.07647708655710783
$ws.Cells.Item(24, 5).Value = 0.06610508080890121
$ws.Cells.Item(24, 6).Value = 1.927709585106072
$ws.Cells.Item(24, 7).Value = 0.002495184121505067
$ws.Cells.Item(24, 9).Value = 1.313477484622688
$ws.Cells.Item(24, 11).Value = 1.40844920406505
$ws.Cells.Item(24, 12).Value = 0.3578700874011815
$ws.Cells.Item(24, 14).Value = 2.340082781948993
$ws.Cells.Item(25, 2).Value = 0.7018354456513123
$ws.Cells.Item(25, 4).Value = 0.07754227193610674
$ws.Cells.Item(25, 5).Value = 0.06342752255451778
$ws.Cells.Item(25, 6).Value = 1.865162919269437
$ws.Cells.Item(25, 7).Value = 0.002503230933536469
$ws.Cells.Item(25, 9).Value = 1.312995808280654
$ws.Cells.Item(25, 11).Value = 1.217075025171141
$ws.Cells.Item(25, 12).Value = 0.3269526318619143
$ws.Cells.Item(25, 14).Value = 2.364579682216572

Write-Output "Applied 216 cell updates"